# The NATMI TPM pipeline was rerun with updated expression values. In the refreshed
# output the "ECs" target-cluster rows no longer qualify, so drop them (originally
# rows 2, 5 and 8) - deleting bottom-to-top keeps the remaining row numbers stable
# while each delete is applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
$ws.Rows(2).Delete()

# Write the recomputed expression / specificity / weight columns (G:T) for the six
# remaining Sending-cluster x Target-cluster pairs (columns A:F already line up
# correctly after the deletions above shifted the surviving rows into place).
$ws.Range("G2").Value2 = 8.898150666666666
$ws.Range("H2").Value2 = 26.694452
$ws.Range("I2").Value2 = 0.3765197173862137
$ws.Range("J2").Value2 = 0.3765197173862137
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 6.072364333333334
$ws.Range("N2").Value2 = 18.217093
$ws.Range("O2").Value2 = 0.4407767221912973
$ws.Range("P2").Value2 = 0.4407767221912974
$ws.Range("Q2").Value2 = 54.03281274089289
$ws.Range("R2").Value2 = 486.295314668036
$ws.Range("S2").Value2 = 0.1659611268698889
$ws.Range("T2").Value2 = 0.1659611268698889
$ws.Range("G3").Value2 = 8.898150666666666
$ws.Range("H3").Value2 = 26.694452
$ws.Range("I3").Value2 = 0.3765197173862137
$ws.Range("J3").Value2 = 0.3765197173862137
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 7.704144333333335
$ws.Range("N3").Value2 = 23.112433
$ws.Range("O3").Value2 = 0.5592232778087027
$ws.Range("P3").Value2 = 0.5592232778087027
$ws.Range("Q3").Value2 = 68.55263703574623
$ws.Range("R3").Value2 = 616.9737333217161
$ws.Range("S3").Value2 = 0.2105585905163248
$ws.Range("T3").Value2 = 0.2105585905163248
$ws.Range("G4").Value2 = 3.344413333333334
$ws.Range("H4").Value2 = 10.03324
$ws.Range("I4").Value2 = 0.1415167724465014
$ws.Range("J4").Value2 = 0.1415167724465015
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 6.072364333333334
$ws.Range("N4").Value2 = 18.217093
$ws.Range("O4").Value2 = 0.4407767221912973
$ws.Range("P4").Value2 = 0.4407767221912974
$ws.Range("Q4").Value2 = 20.30849624125778
$ws.Range("R4").Value2 = 182.77646617132
$ws.Range("S4").Value2 = 0.0623772990940606
$ws.Range("T4").Value2 = 0.06237729909406062
$ws.Range("G5").Value2 = 3.344413333333334
$ws.Range("H5").Value2 = 10.03324
$ws.Range("I5").Value2 = 0.1415167724465014
$ws.Range("J5").Value2 = 0.1415167724465015
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 7.704144333333335
$ws.Range("N5").Value2 = 23.112433
$ws.Range("O5").Value2 = 0.5592232778087027
$ws.Range("P5").Value2 = 0.5592232778087027
$ws.Range("Q5").Value2 = 25.76584303032445
$ws.Range("R5").Value2 = 231.8925872729201
$ws.Range("S5").Value2 = 0.07913947335244084
$ws.Range("T5").Value2 = 0.07913947335244086
$ws.Range("G6").Value2 = 11.39006466666667
$ws.Range("H6").Value2 = 34.170194
$ws.Range("I6").Value2 = 0.4819635101672848
$ws.Range("J6").Value2 = 0.4819635101672848
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 6.072364333333334
$ws.Range("N6").Value2 = 18.217093
$ws.Range("O6").Value2 = 0.4407767221912973
$ws.Range("P6").Value2 = 0.4407767221912974
$ws.Range("Q6").Value2 = 69.16462243622689
$ws.Range("R6").Value2 = 622.4816019260421
$ws.Range("S6").Value2 = 0.2124382962273478
$ws.Range("T6").Value2 = 0.2124382962273478
$ws.Range("G7").Value2 = 11.39006466666667
$ws.Range("H7").Value2 = 34.170194
$ws.Range("I7").Value2 = 0.4819635101672848
$ws.Range("J7").Value2 = 0.4819635101672848
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 7.704144333333335
$ws.Range("N7").Value2 = 23.112433
$ws.Range("O7").Value2 = 0.5592232778087027
$ws.Range("P7").Value2 = 0.5592232778087027
$ws.Range("Q7").Value2 = 87.75070215800024
$ws.Range("R7").Value2 = 789.7563194220021
$ws.Range("S7").Value2 = 0.269525213939937
$ws.Range("T7").Value2 = 0.2695252139399371
